$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 stays empty (matches diff: no row 9 data, gap between row 8 and row 10).
# Add new row 10 - single cell in column A
$ws.Range("A10").Value = "Parametros de tsconfig.json"

# Add new row 11 - columns A and B
$ws.Range("A11").Value = "sourceMap"
$ws.Range("B11").Value = "Permite crear un archivo .map para hacer debug en el explorador"

# Update the selection to match the new active cell location (B12)
$ws.Range("B12").Select()
